{"js": "// New equation text, in reading order, for every non-blank data row of the table\nconst newValues = [\n  \"390\u00f72=195, 0\", \"505\u00f78=63, 1\", \"675\u00f77=96, 3\", \"835\u00f79=92, 7\", \"845\u00f77=120, 5\",\n  \"835\u00f75=167, 0\", \"670\u00f76=111, 4\", \"678\u00f75=135, 3\", \"727\u00f78=90, 7\", \"531\u00f79=59, 0\",\n  \"789\u00f77=112, 5\", \"267\u00f74=66, 3\", \"786\u00f78=98, 2\", \"177\u00f79=19, 6\", \"840\u00f72=420, 0\",\n  \"180\u00f77=25, 5\", \"135\u00f78=16, 7\", \"945\u00f77=135, 0\", \"434\u00f72=217, 0\", \"941\u00f77=134, 3\",\n  \"414\u00f76=69, 0\", \"795\u00f75=159, 0\", \"349\u00f75=69, 4\", \"867\u00f74=216, 3\", \"163\u00f79=18, 1\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nlet k = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  const rowValues = table.values[r];\n  // Spacer rows are fully blank; only the problem rows get new equations\n  const isBlankRow = rowValues.every((cellText) => cellText === \"\");\n  if (isBlankRow) continue;\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[k];\n    k++;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# New equation text, in reading order, for every non-blank data row of the table\n$newValues = @(\n  \"390\u00f72=195, 0\", \"505\u00f78=63, 1\", \"675\u00f77=96, 3\", \"835\u00f79=92, 7\", \"845\u00f77=120, 5\",\n  \"835\u00f75=167, 0\", \"670\u00f76=111, 4\", \"678\u00f75=135, 3\", \"727\u00f78=90, 7\", \"531\u00f79=59, 0\",\n  \"789\u00f77=112, 5\", \"267\u00f74=66, 3\", \"786\u00f78=98, 2\", \"177\u00f79=19, 6\", \"840\u00f72=420, 0\",\n  \"180\u00f77=25, 5\", \"135\u00f78=16, 7\", \"945\u00f77=135, 0\", \"434\u00f72=217, 0\", \"941\u00f77=134, 3\",\n  \"414\u00f76=69, 0\", \"795\u00f75=159, 0\", \"349\u00f75=69, 4\", \"867\u00f74=216, 3\", \"163\u00f79=18, 1\"\n)\n\n$k = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    # Spacer rows only hold the cell-end mark (CR + cell-mark); skip those, only fill problem rows\n    $firstCellText = $tbl.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ([string]::IsNullOrEmpty($firstCellText)) { continue }\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $newValues[$k]\n        $k = $k + 1\n    }\n}"}
